$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Update Crime Complaints table (rows 14-29) ---
$ws.Range("C23").Copy($ws.Range("D14"))
$ws.Range("E23").Copy($ws.Range("E14"))
$ws.Range("H30").Copy($ws.Range("L14"))
$ws.Range("L14").Value = -100
$ws.Range("C23").Copy($ws.Range("C15"))
$ws.Range("H30").Copy($ws.Range("M15"))
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = -63.636363636363
$ws.Range("C16").Value = 2
$ws.Range("C23").Copy($ws.Range("D16"))
$ws.Range("E23").Copy($ws.Range("E16"))
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 16
$ws.Range("K16").Value = -11.111111111111
$ws.Range("L16").Value = 77.777777777777
$ws.Range("M16").Value = -46.666666666666
$ws.Range("N16").Value = -92.558139534883
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 14.814814814814
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 38
$ws.Range("K17").Value = 31.578947368421
$ws.Range("L17").Value = 51.515151515151
$ws.Range("M17").Value = 92.307692307692
$ws.Range("N17").Value = -43.820224719101
$ws.Range("G30").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -66.666666666666
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = -53.333333333333
$ws.Range("L18").Value = -50
$ws.Range("M18").Value = -84.782608695652
$ws.Range("N18").Value = -97.211155378486
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 62.5
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 25
$ws.Range("I19").Value = 61
$ws.Range("J19").Value = 42
$ws.Range("K19").Value = 45.238095238095
$ws.Range("L19").Value = 177.272727272727
$ws.Range("M19").Value = 10.909090909090
$ws.Range("N19").Value = 48.780487804878
$ws.Range("C23").Copy($ws.Range("C20"))
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = -56.521739130434
$ws.Range("N20").Value = -94.444444444444
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = 6.593406593406
$ws.Range("I21").Value = 148
$ws.Range("J21").Value = 130
$ws.Range("K21").Value = 13.846153846153
$ws.Range("L21").Value = 55.789473684210
$ws.Range("M21").Value = -18.232044198895
$ws.Range("N21").Value = -81.289506953223
$ws.Range("F22").Value = 2
$ws.Range("C23").Copy($ws.Range("G22"))
$ws.Range("E23").Copy($ws.Range("H22"))
$ws.Range("I22").Value = 2
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = -66.666666666666
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 5.263157894736
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = 30.263157894736
$ws.Range("I24").Value = 157
$ws.Range("J24").Value = 107
$ws.Range("K24").Value = 46.728971962616
$ws.Range("L24").Value = 7.534246575342
$ws.Range("M24").Value = 28.688524590163
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 40
$ws.Range("H25").Value = 8.108108108108
$ws.Range("I25").Value = 63
$ws.Range("J25").Value = 61
$ws.Range("K25").Value = 3.278688524590
$ws.Range("L25").Value = 70.270270270270
$ws.Range("M25").Value = -18.181818181818
$ws.Range("C23").Copy($ws.Range("C26"))
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 0
$ws.Range("C23").Copy($ws.Range("C27"))
$ws.Range("G30").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("H30").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 6
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = 100
$ws.Range("C23").Copy($ws.Range("D28"))
$ws.Range("E23").Copy($ws.Range("E28"))
$ws.Range("G28").Value = 4
$ws.Range("C23").Copy($ws.Range("D29"))
$ws.Range("E23").Copy($ws.Range("E29"))
$ws.Range("G29").Value = 4
